$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0.0576084665954113, 0.9844096899032593, 0.01507177390158176, 0.998425304889679),
  @(0.008390907198190689, 0.9982026815414429, 0.004987392108887434, 0.9993557929992676),
  @(0.003105229465290904, 0.9991114139556885, 0.001565668499097228, 0.9998568296432495),
  @(0.001660863286815584, 0.9994547367095947, 0.002739168237894773, 0.9994989633560181),
  @(0.001484747277572751, 0.9995759129524231, 0.001401059329509735, 0.9997852444648743),
  @(0.001164978835731745, 0.9997576475143433, 0.001681871013715863, 0.9997852444648743),
  @(0.000867362366989255, 0.9997374415397644, 0.001388045260682702, 0.9998568296432495),
  @(0.001835264731198549, 0.9995557069778442, 0.001665345509536564, 0.9998568296432495),
  @(0.0008348791161552072, 0.9998586177825928, 0.001502249739132822, 0.9997852444648743),
  @(0.0003206375695299357, 0.9999595880508423, 0.001452634576708078, 0.9998568296432495),
  @(0.0003198415506631136, 0.9998586177825928, 0.002353767398744822, 0.9997137188911438),
  @(0.001060804817825556, 0.999798059463501, 0.001451039570383728, 0.9997852444648743),
  @(0.0006114544230513275, 0.9998182654380798, 0.001858364674262702, 0.9998568296432495),
  @(0.0007283990853466094, 0.9998384714126587, 0.002547782380133867, 0.9997852444648743),
  @(0.0001959305809577927, 0.9999394416809082, 0.002105168765410781, 0.9997852444648743),
  @(0.0005859470693394542, 0.9999192357063293, 0.0008411071030423045, 0.9997852444648743),
  @(0.0004435898154042661, 0.9998586177825928, 0.001463642111048102, 0.9998568296432495),
  @(0.0003710883029270917, 0.9998788237571716, 0.001283077173866332, 0.9997852444648743),
  @(0.0002178992726840079, 0.9999394416809082, 0.002492540050297976, 0.9998568296432495),
  @(0.00022945660748519, 0.9998990297317505, 0.002242187736555934, 0.9997852444648743),
  @(0.0003285407146904618, 0.9999192357063293, 0.001912795123644173, 0.9997852444648743),
  @(0.00008283296483568847, 0.9999595880508423, 0.002252142177894711, 0.9998568296432495),
  @(0.00002261812733195256, 1, 0.001655915868468583, 0.9997852444648743),
  @(0.001296588568948209, 0.9997172951698303, 0.0003999467589892447, 0.9999284148216248),
  @(0.00007327401544898748, 0.9999797940254211, 0.001411285367794335, 0.9997852444648743),
  @(0.000190761114936322, 0.9999797940254211, 0.002003897912800312, 0.9997852444648743),
  @(0.0005537224933505058, 0.9998990297317505, 0.0002932380884885788, 0.9998568296432495),
  @(0.00002671337279025465, 1, 0.0004444809164851904, 0.9997852444648743),
  @(0.0001308671344304457, 0.9999595880508423, 0.0004463412333279848, 0.9999284148216248),
  @(0.00001320227238466032, 1, 0.0006775259971618652, 0.9998568296432495),
  @(0.0005253042909316719, 0.9998384714126587, 0.0001497242919867858, 0.9999284148216248),
  @(0.0005554304807446897, 0.9998990297317505, 0.0001786122738849372, 0.9998568296432495),
  @(0.00003182193904649466, 0.9999797940254211, 0.001665437826886773, 0.9998568296432495),
  @(0.0002613798424135894, 0.9999595880508423, 0.001240541110746562, 0.9998568296432495),
  @(0.0003156974271405488, 0.9999595880508423, 0.0008084229193627834, 0.9998568296432495),
  @(0.0002196431596530601, 0.9999595880508423, 0.001076194457709789, 0.9998568296432495),
  @(0.0002906432200688869, 0.9999394416809082, 0.0006955883000046015, 0.9997852444648743),
  @(0.00001505835825810209, 1, 0.0006742156110703945, 0.9997852444648743),
  @(0.000001918103180287289, 1, 0.0007515454781241715, 0.9997852444648743),
  @(0.0003434353566262871, 0.9999394416809082, 0.002784417476505041, 0.9997852444648743),
  @(0.0006921313470229506, 0.9997778534889221, 0.001658034627325833, 0.9998568296432495),
  @(0.0004707902844529599, 0.9998384714126587, 0.0009196344180963933, 0.9998568296432495),
  @(0.0004515611799433827, 0.9999595880508423, 0.002744687488302588, 0.9998568296432495),
  @(0.0002983151935040951, 0.9999394416809082, 0.002172876382246614, 0.9998568296432495),
  @(0.0002042712731054053, 0.9999797940254211, 0.002478931332007051, 0.9998568296432495),
  @(0.000006489568477263674, 1, 0.002649120520800352, 0.9998568296432495),
  @(0.0002376199554419145, 0.9998788237571716, 0.002420897828415036, 0.9998568296432495),
  @(0.000182044692337513, 0.9999192357063293, 0.00004250204074196517, 1),
  @(0.0006724332342855632, 0.9998586177825928, 0.0009184758528135717, 0.9998568296432495),
  @(0.00002023386332439259, 1, 0.001325072953477502, 0.9998568296432495)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt 4; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}

$wb.Save()
